# Apply the diff: reset several "Features" rows' precision/recall/f1/sem
# text metrics to "0,000" and adjust their tp/fp/fn integer counters, then
# refresh two rows' count/pct summaries, and finally zero the
# "Global Metrics" sheet's F1/SEM/W F1/W SEM row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Features")
$gm = $wb.Worksheets.Item("Global Metrics")

function Set-TextCell($sheet, $addr, $text) {
    # Leading apostrophe forces Excel to store the value as literal text
    # (quote-prefix) instead of re-parsing "0,000" style strings as numbers.
    $sheet.Range($addr).Value = "'" + $text
}

# ---- Features sheet ----

# Row 2
Set-TextCell $ws "B2" "0,000"
Set-TextCell $ws "C2" "0,000"
Set-TextCell $ws "D2" "0,000"
Set-TextCell $ws "E2" "0,000"
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 30

# Row 3
Set-TextCell $ws "B3" "0,000"
Set-TextCell $ws "C3" "0,000"
Set-TextCell $ws "D3" "0,000"
Set-TextCell $ws "E3" "0,000"
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 28

# Row 4
Set-TextCell $ws "B4" "0,000"
Set-TextCell $ws "C4" "0,000"
Set-TextCell $ws "D4" "0,000"
Set-TextCell $ws "E4" "0,000"
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 31

# Row 5
Set-TextCell $ws "B5" "0,000"
Set-TextCell $ws "C5" "0,000"
Set-TextCell $ws "D5" "0,000"
Set-TextCell $ws "E5" "0,000"
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 29

# Row 6
$ws.Range("I6").Value = 0

# Row 7
Set-TextCell $ws "B7" "0,000"
Set-TextCell $ws "C7" "0,000"
Set-TextCell $ws "D7" "0,000"
Set-TextCell $ws "E7" "0,000"
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 27

# Row 8
Set-TextCell $ws "B8" "0,000"
Set-TextCell $ws "C8" "0,000"
Set-TextCell $ws "D8" "0,000"
Set-TextCell $ws "E8" "0,000"
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 13

# Row 11
Set-TextCell $ws "B11" "0,000"
Set-TextCell $ws "C11" "0,000"
Set-TextCell $ws "D11" "0,000"
Set-TextCell $ws "E11" "0,000"
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 11

# Row 12
Set-TextCell $ws "B12" "0,000"
Set-TextCell $ws "C12" "0,000"
Set-TextCell $ws "D12" "0,000"
Set-TextCell $ws "E12" "0,000"
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 11

# Row 13
Set-TextCell $ws "B13" "0,000"
Set-TextCell $ws "C13" "0,000"
Set-TextCell $ws "D13" "0,000"
Set-TextCell $ws "E13" "0,000"
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 12

# Row 14
Set-TextCell $ws "B14" "0,000"
Set-TextCell $ws "C14" "0,000"
Set-TextCell $ws "D14" "0,000"
Set-TextCell $ws "E14" "0,000"
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 11

# Row 22
Set-TextCell $ws "F22" "10,000"
Set-TextCell $ws "G22" "0,18"
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 10

# Row 38
Set-TextCell $ws "F38" "3,000"
Set-TextCell $ws "G38" "0,05"
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 3

# ---- Global Metrics sheet ----

Set-TextCell $gm "B2" "0,000"
Set-TextCell $gm "C2" "0,000"
Set-TextCell $gm "D2" "0,000"
Set-TextCell $gm "E2" "0,000"
